$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the CasesTab query (B2): remove Cohort-related lines ---
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in [ 'Bluetick Hound', 'Norfolk Terrier', 'Scottish Terrier']and diag.disease_term in ['Bladder Cancer'] and demo.sex in ['Male', 'Female'] and demo.neutered_indicator IN ['No']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@
$ws.Range("B2").Value = $casesQuery

# --- Update the SamplesTab query (B3): now uses the former index-12 text (unchanged content) ---
$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
	WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in [ 'Bluetick Hound', 'Norfolk Terrier', 'Scottish Terrier']and diag.disease_term in ['Bladder Cancer'] and demo.sex in ['Male', 'Female'] and demo.neutered_indicator IN ['No']
 WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@
$ws.Range("B3").Value = $samplesQuery

# --- Update the FilesTab query (B4): now uses the former index-13 text (unchanged content) ---
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in [ 'Bluetick Hound', 'Norfolk Terrier', 'Scottish Terrier']and diag.disease_term in ['Bladder Cancer'] and demo.sex in ['Male', 'Female'] and demo.neutered_indicator IN ['No']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(f.file_type, '') AS `File Type`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `File Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@
$ws.Range("B4").Value = $filesQuery

# --- Row heights reflect the re-wrapped text after the edits ---
$ws.Rows.Item(2).RowHeight = 290
$ws.Rows.Item(3).RowHeight = 261
$ws.Rows.Item(4).RowHeight = 261

# --- Update the view: scroll/select to match the new focal cell ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B2").Select()

